$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without letting Excel
# reinterpret numeric-looking text (e.g. "50.4900") as a number. Only
# cells whose current display format is numeric need a temporary switch
# to Text ("@") while the value is written; cells already formatted as
# Text can be set directly. Re-reading/re-applying NumberFormat only when
# it actually needs to change avoids Excel re-resolving the cell style to
# a different (but numerically-equivalent) style record, which would
# otherwise quietly change alignment/reading-order.
function Set-TextValue($range, $text) {
    $fmt = $range.NumberFormat
    if ($fmt -ne "@") {
        $range.NumberFormat = "@"
        $range.Value = $text
        $range.NumberFormat = $fmt
    } else {
        $range.Value = $text
    }
}

# Row 26 - "VOLTAREN 75MG/3ML 3 AMP." : balance 4:1 -> 3:2, sell price 16.8300 -> 50.4900,
# transactions 0:1 -> 0:3
Set-TextValue $ws.Range("H26") "3:2"
Set-TextValue $ws.Range("P26") "50.4900"
Set-TextValue $ws.Range("Q26") "0:3"

# Row 30 - "سرنجات 3 سم" : sell price 4.0000 -> 10.0000, transactions 2:0 -> 5:0
Set-TextValue $ws.Range("P30") "10.0000"
Set-TextValue $ws.Range("Q30") "5:0"

# Row 37 - total of the P column
$ws.Range("P37").Value = 1066.0999999999999

# Row 38 - footer timestamp updated to reflect the new export time
Set-TextValue $ws.Range("A38") "Wednesday, 20 August, 2025 5:05 PM"
